# Fix the "decimal point replaced by comma" data entry bug.
# The cells C4:F8 (excluding the "NA" text cell) held text values of
# numbers (stored as shared strings). Re-enter them as real numeric
# values, correcting the one cell (D5) where a comma had replaced the
# decimal point (23,34 -> 23.34 instead of the erroneous 23.24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = 6.83
$ws.Range("D4").Value = 6.8
$ws.Range("F4").Value = 6.71

$ws.Range("C5").Value = 23.34
$ws.Range("D5").Value = 23.34
$ws.Range("E5").Value = 23.2
$ws.Range("F5").Value = 23.68

$ws.Range("C6").Value = 19.35
$ws.Range("D6").Value = 19.38
$ws.Range("E6").Value = 19.58
$ws.Range("F6").Value = 19.68

$ws.Range("C7").Value = 3.84
$ws.Range("D7").Value = 3.84
$ws.Range("E7").Value = 3.54
$ws.Range("F7").Value = 3.74

$ws.Range("C8").Value = 19.38
$ws.Range("D8").Value = 19.58
$ws.Range("E8").Value = 19.35
$ws.Range("F8").Value = 19.69

# Restore the selection to C4, matching the saved workbook state.
$ws.Range("C4").Select()

# Page setup as saved by Excel (paper size / orientation) for this sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
